# Weekly fruit/vegetable price update: rotate data among rows 2, 3 and 4
# New row 2 = old row 4 data, new row 3 = old row 2 data, new row 4 = old row 3 data
# (only columns D, H, K, L, M, N, O, P change; other columns are identical across
# these rows so they don't need to be touched)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "old" values before any write (use Value2 - plain COM Value
# reflects oddly in this host)
$oldD2 = $ws.Range("D2").Value2
$oldH2 = $ws.Range("H2").Value2
$oldK2 = $ws.Range("K2").Value2
$oldL2 = $ws.Range("L2").Value2
$oldM2 = $ws.Range("M2").Value2
$oldN2 = $ws.Range("N2").Value2
$oldO2 = $ws.Range("O2").Value2
$oldP2 = $ws.Range("P2").Value2

$oldD3 = $ws.Range("D3").Value2
$oldH3 = $ws.Range("H3").Value2
$oldK3 = $ws.Range("K3").Value2
$oldL3 = $ws.Range("L3").Value2
$oldM3 = $ws.Range("M3").Value2
$oldN3 = $ws.Range("N3").Value2
$oldO3 = $ws.Range("O3").Value2
$oldP3 = $ws.Range("P3").Value2

$oldD4 = $ws.Range("D4").Value2
$oldH4 = $ws.Range("H4").Value2
$oldK4 = $ws.Range("K4").Value2
$oldL4 = $ws.Range("L4").Value2
$oldM4 = $ws.Range("M4").Value2
$oldN4 = $ws.Range("N4").Value2
$oldO4 = $ws.Range("O4").Value2
$oldP4 = $ws.Range("P4").Value2

# Row 2 <- old Row 4
$ws.Range("D2").Value2 = $oldD4
$ws.Range("H2").Value2 = $oldH4
$ws.Range("K2").Value2 = $oldK4
$ws.Range("L2").Value2 = $oldL4
$ws.Range("M2").Value2 = $oldM4
$ws.Range("N2").Value2 = $oldN4
$ws.Range("O2").Value2 = $oldO4
$ws.Range("P2").Value2 = $oldP4

# Row 3 <- old Row 2
$ws.Range("D3").Value2 = $oldD2
$ws.Range("H3").Value2 = $oldH2
$ws.Range("K3").Value2 = $oldK2
$ws.Range("L3").Value2 = $oldL2
$ws.Range("M3").Value2 = $oldM2
$ws.Range("N3").Value2 = $oldN2
$ws.Range("O3").Value2 = $oldO2
$ws.Range("P3").Value2 = $oldP2

# Row 4 <- old Row 3
$ws.Range("D4").Value2 = $oldD3
$ws.Range("H4").Value2 = $oldH3
$ws.Range("K4").Value2 = $oldK3
$ws.Range("L4").Value2 = $oldL3
$ws.Range("M4").Value2 = $oldM3
$ws.Range("N4").Value2 = $oldN3
$ws.Range("O4").Value2 = $oldO3
$ws.Range("P4").Value2 = $oldP3
